$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the description for row 3: "Research for inspiration" -> "Research for inspiration & plugins"
$ws.Range("A3").Value = "Research for inspiration & plugins"

# Update hours for row 3: 0.5 -> 1
$ws.Range("C3").Value = 1

# Update the active selection to F7
$ws.Range("F7").Select()
